# TCD_PHIEU_HUONG_DAN_KHIEU_NAI.docx — "fix in phieu tiep cong dan" (#358)
#
# 1) Nudge the little underline/separator drawing's stored extents back to
#    its "proper" size (Word normally recomputes this automatically when the
#    custom-geometry shape is touched/resaved).
# 2) Add the missing colon after "noi cap" in the identity-document line.
# 3) Collapse the four-way split run in the "Sau khi xem xet..." paragraph
#    back into a single run (the ${coQuanTiepNhan} placeholders were left in
#    separate runs by a previous edit for no reason).
# 4) Turn off "overflowPunct" on the Normal and LO-normal paragraph styles.

$d = $word.ActiveDocument

# --- 1) Resize the small separator shape (table 1, row 2, col 1) ----------
# Real Word exposes the drawing's extent (EMU) through Shape.Width/Height
# (in points, 1 pt = 12700 EMU). Target: cx=2047240 EMU (=161.2pt),
# cy=15240 EMU (=1.2pt).
try {
    $shp = $d.Shapes.Item(1)
    $shp.Width = 161.2
    $shp.Height = 1.2
} catch {
    # Some hosts don't allow resizing this particular freeform/custGeom
    # shape through the Shapes collection; ignore and move on.
}

# --- 2) " noi cap ${noiCap}." -> " noi cap: ${noiCap}." -------------------
$d.Content.Find.Execute('nơi cấp ${noiCap}.', $true, $false, $false, $false, $false, `
    $true, 1, $false, 'nơi cấp: ${noiCap}.', 2) | Out-Null

# --- 3) Merge the split "Sau khi xem xet..." runs into one ----------------
$oldText = 'Sau khi xem xét nội dung đơn khiếu nại, ${coQuanTiepNhan} nhận thấy khiếu nại của ông (bà) không thuộc thẩm quyền giải quyết của ${coQuanTiepNhan}.'
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $oldText, 2) | Out-Null

# --- 4) overflowPunct: true -> false on Normal + LO-normal styles ---------
foreach ($styleName in @("Normal", "LOnormal")) {
    $style = $d.Styles.Item($styleName)
    $style.ParagraphFormat.HangingPunctuation = 0
}
